# Update NATMI LR-pair TPM-derived metrics for Rspo3-Fzd8 with recalculated
# values from the new TPM matrix. Only the numeric metric columns (G:T) are
# affected; identifier columns (A:F) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.003420333333333333
$ws.Range("H2").Value = 0.010261
$ws.Range("I2").Value = 0.003549653112303053
$ws.Range("J2").Value = 0.003549653112303053
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.235341333333333
$ws.Range("N2").Value = 9.706024
$ws.Range("O2").Value = 0.2153734454473681
$ws.Range("P2").Value = 0.2153734454473681
$ws.Range("Q2").Value = 0.01106594580711111
$ws.Range("R2").Value = 0.09959351226399998
$ws.Range("S2").Value = 0.000764501020939682
$ws.Range("T2").Value = 0.0007645010209396821

# Row 3
$ws.Range("G3").Value = 0.003420333333333333
$ws.Range("H3").Value = 0.010261
$ws.Range("I3").Value = 0.003549653112303053
$ws.Range("J3").Value = 0.003549653112303053
$ws.Range("O3").Value = 0.4841904166376352
$ws.Range("P3").Value = 0.4841904166376352
$ws.Range("Q3").Value = 0.02487783440388889
$ws.Range("R3").Value = 0.223900509635
$ws.Range("S3").Value = 0.001718708019365094
$ws.Range("T3").Value = 0.001718708019365094

# Row 4
$ws.Range("G4").Value = 0.003420333333333333
$ws.Range("H4").Value = 0.010261
$ws.Range("I4").Value = 0.003549653112303053
$ws.Range("J4").Value = 0.003549653112303053
$ws.Range("O4").Value = 0.3004361379149967
$ws.Range("P4").Value = 0.3004361379149967
$ws.Range("Q4").Value = 0.01543648992455555
$ws.Range("R4").Value = 0.138928409321
$ws.Range("S4").Value = 0.001066444071998278
$ws.Range("T4").Value = 0.001066444071998278

# Row 5
$ws.Range("I5").Value = 0.3907064193682856
$ws.Range("J5").Value = 0.3907064193682855
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.235341333333333
$ws.Range("N5").Value = 9.706024
$ws.Range("O5").Value = 0.2153734454473681
$ws.Range("P5").Value = 0.2153734454473681
$ws.Range("Q5").Value = 1.218016500889778
$ws.Range("R5").Value = 10.962148508008
$ws.Range("S5").Value = 0.08414778769775197
$ws.Range("T5").Value = 0.08414778769775197

# Row 6
$ws.Range("I6").Value = 0.3907064193682856
$ws.Range("J6").Value = 0.3907064193682855
$ws.Range("O6").Value = 0.4841904166376352
$ws.Range("P6").Value = 0.4841904166376352
$ws.Range("S6").Value = 0.1891763039769288
$ws.Range("T6").Value = 0.1891763039769288

# Row 7
$ws.Range("I7").Value = 0.3907064193682856
$ws.Range("J7").Value = 0.3907064193682855
$ws.Range("O7").Value = 0.3004361379149967
$ws.Range("P7").Value = 0.3004361379149967
$ws.Range("S7").Value = 0.1173823276936048
$ws.Range("T7").Value = 0.1173823276936048

# Row 8
$ws.Range("I8").Value = 0.6057439275194114
$ws.Range("J8").Value = 0.6057439275194113
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.235341333333333
$ws.Range("N8").Value = 9.706024
$ws.Range("O8").Value = 0.2153734454473681
$ws.Range("P8").Value = 0.2153734454473681
$ws.Range("Q8").Value = 1.888390009627555
$ws.Range("R8").Value = 16.995510086648
$ws.Range("S8").Value = 0.1304611567286764
$ws.Range("T8").Value = 0.1304611567286764

# Row 9
$ws.Range("I9").Value = 0.6057439275194114
$ws.Range("J9").Value = 0.6057439275194113
$ws.Range("O9").Value = 0.4841904166376352
$ws.Range("P9").Value = 0.4841904166376352
$ws.Range("S9").Value = 0.2932954046413414
$ws.Range("T9").Value = 0.2932954046413413

# Row 10
$ws.Range("I10").Value = 0.6057439275194114
$ws.Range("J10").Value = 0.6057439275194113
$ws.Range("O10").Value = 0.3004361379149967
$ws.Range("P10").Value = 0.3004361379149967
$ws.Range("S10").Value = 0.1819873661493937
$ws.Range("T10").Value = 0.1819873661493936

